$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: DAMSLTag sv -> ba, DialogAct Statement-opinion -> Appreciation
$ws.Range("I6").Value = "ba"
$ws.Range("J6").Value = "Appreciation"

# Row 30: DAMSLTag aa -> sd, DialogAct Agree/Accept -> Statement-non-opinion
$ws.Range("I30").Value = "sd"
$ws.Range("J30").Value = "Statement-non-opinion"

# Row 33: DAMSLTag sd -> sv, DialogAct Statement-non-opinion -> Statement-opinion
$ws.Range("I33").Value = "sv"
$ws.Range("J33").Value = "Statement-opinion"
